$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.337.84'
$ws.Range('E2').Value = '  -0.17%  '
$ws.Range('D3').Value = '1.839.86'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '239.14'
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('D6').Value = '0.6281'
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').Value = '0.07428'
$ws.Range('E8').Value = '  -1.05%  '
$ws.Range('D9').Value = '24.98'
$ws.Range('E9').Value = '  +2.19%  '
$ws.Range('D10').Value = '0.2889'
$ws.Range('E10').Value = '  -0.68%  '
$ws.Range('D11').Value = '0.07730'
$ws.Range('E11').Value = '  +0.25%  '
$ws.Range('D12').Value = '1.804.62'
$ws.Range('E12').Value = '  -2.31%  '
$ws.Range('D13').Value = '4.955'
$ws.Range('E13').Value = '  -1.05%  '
$ws.Range('D14').Value = '0.6736'
$ws.Range('E14').Value = '  -1.08%  '
$ws.Range('D15').Value = '0.00001021'
$ws.Range('E15').Value = '  -0.94%  '
$ws.Range('D16').Value = '81.55'
$ws.Range('E16').Value = '  -0.83%  '
$ws.Range('D17').Value = '6.211'
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('D18').Value = '29.337.96'
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('D19').Value = '228.77'
$ws.Range('E19').Value = '  -0.36%  '
$ws.Range('E20').Value = '  -0.60%  '
$ws.Range('E21').Value = '  +0.15%  '
$ws.Range('D22').Value = '7.330'
$ws.Range('E22').Value = '  -1.67%  '
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').Value = '158.10'
$ws.Range('E24').Value = '  -0.43%  '
$ws.Range('D25').Value = '8.457'
$ws.Range('E25').Value = '  +0.46%  '
$ws.Range('D26').Value = '0.1344'
$ws.Range('E26').Value = '  -2.64%  '
$ws.Range('D27').Value = '17.35'
$ws.Range('E27').Value = '  -1.27%  '
$ws.Range('D28').Value = '0.07405'
$ws.Range('E28').Value = '  +15.33%  '
$ws.Range('D29').Value = '1.461'
$ws.Range('E29').Value = '  +5.27%  '
$ws.Range('D30').Value = '1.476'
$ws.Range('E30').Value = '  +0.21%  '
$ws.Range('D31').Value = '4.026'
$ws.Range('E31').Value = '  -1.70%  '
$ws.Range('D32').Value = '4.034'
$ws.Range('E32').Value = '  -0.68%  '
$ws.Range('D33').Value = '1.818'
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('D34').Value = '1.137'
$ws.Range('E34').Value = '  -0.48%  '
$ws.Range('D35').Value = '0.6935'
$ws.Range('E35').Value = '  -0.86%  '
$ws.Range('D36').Value = '2.577'
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('E37').Value = '  +0.63%  '
$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38').Value = '6.884'
$ws.Range('E38').Value = '  +4.55%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '2.806'
$ws.Range('E39').Value = '  -1.00%  '
$ws.Range('D40').Value = '1.231.33'
$ws.Range('E40').Value = '  -2.32%  '
$ws.Range('D41').Value = '0.9289'
$ws.Range('E41').Value = '  +2.11%  '
$ws.Range('D42').Value = '1.000'
$ws.Range('E42').Value = '  +0.19%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = '100.59'
$ws.Range('E43').Value = '  -0.77%  '
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '1.971.55'
$ws.Range('E44').Value = '  -1.76%  '
$ws.Range('D45').Value = '65.19'
$ws.Range('E45').Value = '  -1.62%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.00000000120'
$ws.Range('E46').Value = '  +0.91%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').Value = '1.701'
$ws.Range('E47').Value = '  -0.73%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').Value = '6.925'
$ws.Range('E48').Value = '  -1.98%  '
$ws.Range('D49').Value = '0.1137'
$ws.Range('E49').Value = '  -3.42%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '8.861'
$ws.Range('E50').Value = '  -2.27%  '
$ws.Range('B51').Value = 'TheSandbox'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D51').Value = '0.3899'
$ws.Range('E51').Value = '  -1.13%  '
